# "RMI files through Dec 1"
# Update the Government Revenue Accounting weights for "carbon tax revenue"
# (row 8 on the "Set Values Here" sheet): split the "Deficit Spending" and
# "Payroll Taxes" weights so the row reads 0, 5, 5, 5, 0 instead of 0, 0, 5, 0, 0.

$wb = $excel.ActiveWorkbook

$setValues = $wb.Worksheets.Item("Set Values Here")
$setValues.Activate()

$setValues.Range("C8").Value = 5
$setValues.Range("E8").Value = 5

$setValues.Range("C9").Select()

$carbonTax = $wb.Worksheets.Item("GRA-carbontax")
$carbonTax.Activate()
$carbonTax.Range("B5").Select()

$about = $wb.Worksheets.Item("About")
$about.Activate()
